$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 19
$ws.Range("F4").Value = 585
$ws.Range("F6").Value = 27
$ws.Range("F7").Value = 1925
$ws.Range("F8").Value = 5272
$ws.Range("F9").Value = 1468
$ws.Range("F11").Value = 3070
$ws.Range("F13").Value = 35
$ws.Range("F14").Value = 1277
$ws.Range("F15").Value = 4202
$ws.Range("F16").Value = 1001
$ws.Range("F18").Value = 1646
$ws.Range("F19").Value = 2591
$ws.Range("F20").Value = 27
$ws.Range("F22").Value = 125
$ws.Range("F23").Value = 142
$ws.Range("F24").Value = 958
$ws.Range("F25").Value = 286
$ws.Range("F27").Value = 75
$ws.Range("F29").Value = 1068
$ws.Range("F30").Value = 355
$ws.Range("F31").Value = 29
$ws.Range("F32").Value = 117
$ws.Range("F34").Value = 231
$ws.Range("F35").Value = 1632
$ws.Range("F36").Value = 2151
$ws.Range("F37").Value = 1004
$ws.Range("F38").Value = 33
$ws.Range("F39").Value = 244
$ws.Range("F40").Value = 596
$ws.Range("F41").Value = 268
$ws.Range("F43").Value = 647
$ws.Range("F44").Value = 391
$ws.Range("F45").Value = 303
$ws.Range("F46").Value = 204
$ws.Range("F47").Value = 128

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 725

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 725
$ws.Range("F4").Value = 19
$ws.Range("F6").Value = 585
$ws.Range("F7").Value = 27
$ws.Range("F8").Value = 1925
$ws.Range("F9").Value = 5272
$ws.Range("F10").Value = 1468
$ws.Range("F12").Value = 12
$ws.Range("F13").Value = 3070
$ws.Range("F14").Value = 35
$ws.Range("F15").Value = 1277
$ws.Range("F16").Value = 4202
$ws.Range("F17").Value = 1001
$ws.Range("F18").Value = 1646
$ws.Range("F20").Value = 2591
$ws.Range("F22").Value = 27
$ws.Range("F25").Value = 142
$ws.Range("F27").Value = 958
$ws.Range("F28").Value = 286
$ws.Range("F29").Value = 75
$ws.Range("F32").Value = 1068
$ws.Range("F33").Value = 355
$ws.Range("F34").Value = 29
$ws.Range("F36").Value = 1632
$ws.Range("F37").Value = 2151
$ws.Range("F38").Value = 1004
$ws.Range("F39").Value = 33
$ws.Range("F41").Value = 244
$ws.Range("F42").Value = 596
$ws.Range("F43").Value = 268
$ws.Range("F44").Value = 647
$ws.Range("F45").Value = 391
$ws.Range("F46").Value = 303
$ws.Range("F47").Value = 204
$ws.Range("F48").Value = 128

